# Daily attendance processing - reorder "Recorded By" names in column G.
# For every cell in column G that contains a comma-separated list of names,
# rotate the list left by one position (move the first entry to the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value = $rotated
        }
    }
}
